# Refresh the cryptocurrency price/volume table with the latest scraped values.
# Generated from the upstream GitHub Actions data-refresh commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text happens to parse as a plain number (e.g. "13.00", "1.00").
# Excel's COM layer auto-converts such strings to numeric cells, which would both
# change the stored cell type away from text and normalize the displayed digits
# ("13.00" -> 13, "603.46" -> 603.46000000000004). Forcing a text number format
# while assigning the value keeps it a plain text string, then resetting the style
# back to Normal removes the temporary formatting so the cell looks untouched.
$textNumberCells = @{
    D5 = "603.46"
    D6 = "197.19"
    D7 = "0.628"
    D9 = "0.206"
    D11 = "53.60"
    D13 = "9.57"
    D15 = "599.15"
    D16 = "13.00"
    D19 = "19.05"
    D21 = "0.996"
    D22 = "17.91"
    D23 = "5.17"
    D24 = "101.82"
    D26 = "3.01"
    D27 = "10.75"
    D28 = "9.65"
    D29 = "33.88"
    D34 = "63.44"
    D37 = "538.19"
    D38 = "3.11"
    D39 = "1.00"
    D40 = "36.97"
    D42 = "3.54"
    D48 = "8.59"
    D51 = "1.30"
}

foreach ($addr in $textNumberCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $textNumberCells[$addr]
    $cell.Style = "Normal"
}

# Remaining cells: plain text updates (prices with thousands separators, and the
# padded percentage strings in column E) that Excel already stores as text as-is.
$textCells = @{
    D2 = "70.505.52"
    E2 = "  +0.69%  "
    D3 = "3.623.77"
    E3 = "  +2.33%  "
    E5 = "  -0.07%  "
    E6 = "  -0.02%  "
    E7 = "  +0.30%  "
    E9 = "  -1.19%  "
    E10 = "  -0.87%  "
    E11 = "  -0.85%  "
    E12 = "  +0.51%  "
    E13 = "  +0.42%  "
    D14 = "4.196.67"
    E14 = "  +2.27%  "
    E15 = "  -0.50%  "
    E16 = "  +1.37%  "
    D17 = "70.623.81"
    E17 = "  +0.67%  "
    D18 = "3.632.38"
    E18 = "  +2.50%  "
    E19 = "  -1.12%  "
    E20 = "  +1.65%  "
    E21 = "  +0.15%  "
    E22 = "  +0.14%  "
    E23 = "  -1.62%  "
    E24 = "  -0.47%  "
    E25 = "  +0.22%  "
    E26 = "  -3.81%  "
    E27 = "  -1.85%  "
    E28 = "  +0.68%  "
    E29 = "  +0.67%  "
    E30 = "  +7.65%  "
    E31 = "  +1.41%  "
    E32 = "  -2.77%  "
    E33 = "  +2.82%  "
    D35 = "0.0₃0887"
    E35 = "  +5.77%  "
    D36 = "3.893.75"
    E36 = "  +2.97%  "
    E37 = "  +8.66%  "
    E38 = "  +1.12%  "
    E39 = "  -0.07%  "
    E40 = "  +0.81%  "
    E41 = "  -1.00%  "
    E42 = "  -2.94%  "
    E43 = "  +0.14%  "
    E44 = "  +0.60%  "
    E45 = "  +3.71%  "
    E46 = "  +0.94%  "
    E47 = "  +0.00%  "
    E48 = "  -0.90%  "
    E50 = "  +0.68%  "
    E51 = "  +1.30%  "
}

foreach ($addr in $textCells.Keys) {
    $ws.Range($addr).Value = $textCells[$addr]
}

